$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated Price values are plain numeric strings (e.g. "243.38").
# Column D stores prices as TEXT, so pre-format those specific cells as
# Text before writing, otherwise Excel auto-converts them to numbers.
$textCells = @("D4", "D5", "D6", "D8", "D9", "D10", "D13", "D14", "D15", "D16", "D18", "D21", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values (Coin / Link / Price / Volume(1h))
$ws.Range('D2').Value = '29.575.71'
$ws.Range('E2').Value = '  -0.46%  '
$ws.Range('D3').Value = '1.852.69'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '243.38'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('D6').Value = '0.6391'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.07588'
$ws.Range('E8').Value = '  +1.55%  '
$ws.Range('D9').Value = '0.3006'
$ws.Range('E9').Value = '  +0.74%  '
$ws.Range('D10').Value = '24.33'
$ws.Range('E10').Value = '  +0.27%  '
$ws.Range('E11').Value = '  +0.36%  '
$ws.Range('D12').Value = '1.865.12'
$ws.Range('E12').Value = '  +0.46%  '
$ws.Range('D13').Value = '5.050'
$ws.Range('E13').Value = '  +0.12%  '
$ws.Range('D14').Value = '0.6895'
$ws.Range('E14').Value = '  +0.39%  '
$ws.Range('D15').Value = '84.13'
$ws.Range('E15').Value = '  +0.51%  '
$ws.Range('D16').Value = '0.000009702'
$ws.Range('E16').Value = '  +2.02%  '
$ws.Range('D17').Value = '2.113.93'
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('D18').Value = '6.294'
$ws.Range('E18').Value = '  +3.98%  '
$ws.Range('D19').Value = '29.610.24'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range('D21').Value = '12.58'
$ws.Range('E21').Value = '  -0.26%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').Value = '7.647'
$ws.Range('E23').Value = '  +3.37%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').Value = '156.91'
$ws.Range('E25').Value = '  -0.76%  '
$ws.Range('D26').Value = '0.1399'
$ws.Range('E26').Value = '  -1.19%  '
$ws.Range('D27').Value = '8.505'
$ws.Range('E27').Value = '  +0.17%  '
$ws.Range('D28').Value = '17.78'
$ws.Range('E28').Value = '  -0.58%  '
$ws.Range('D29').Value = '1.488'
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('D30').Value = '0.05928'
$ws.Range('E30').Value = '  -4.98%  '
$ws.Range('E31').Value = '  +1.06%  '
$ws.Range('D32').Value = '4.138'
$ws.Range('E32').Value = '  -0.20%  '
$ws.Range('D33').Value = '4.085'
$ws.Range('E33').Value = '  -0.29%  '
$ws.Range('D34').Value = '1.907'
$ws.Range('E34').Value = '  +1.40%  '
$ws.Range('D35').Value = '1.183'
$ws.Range('E35').Value = '  +0.81%  '
$ws.Range('D36').Value = '0.7228'
$ws.Range('E36').Value = '  -0.59%  '
$ws.Range('D37').Value = '2.600'
$ws.Range('E37').Value = '  -0.24%  '
$ws.Range('D38').Value = '2.803'
$ws.Range('E38').Value = '  -1.46%  '
$ws.Range('D39').Value = '1.220.72'
$ws.Range('E39').Value = '  +1.44%  '
$ws.Range('D40').Value = '0.01778'
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('D41').Value = '0.9143'
$ws.Range('E41').Value = '  -1.09%  '
$ws.Range('D42').Value = '6.136'
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('D44').Value = '2.018.74'
$ws.Range('E44').Value = '  -0.52%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = '101.93'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '67.44'
$ws.Range('E46').Value = '  +2.19%  '
$ws.Range('D47').Value = '7.443'
$ws.Range('E47').Value = '  +11.37%  '
$ws.Range('D48').Value = '0.4063'
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('E49').Value = '  -3.23%  '
$ws.Range('D50').Value = '9.152'
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('D51').Value = '1.678'
$ws.Range('E51').Value = '  +1.94%  '

# Restore the default cell style on the cells we temporarily reformatted,
# so only the values (not the formatting) differ from the original file.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
